$d = $word.ActiveDocument

# The "Địa chỉ: " bullet item in Bên A's info block had a second run
# holding the unresolved merge placeholder "vnpt.SiteAddress" right
# after it. Find that placeholder run's text and delete it outright
# (not just clear its text), so the whole <w:r> collapses away and the
# paragraph reads just "Địa chỉ: ".
$range = $d.Content
$found = $range.Find.Execute("vnpt.SiteAddress", $true, $false, $false, $false, $false, $true, 1, $false, "")
if ($found) {
    $range.Delete()
}
